# Add season record columns (Wins, Losses, Ties) to the PHI 2009 roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the last existing header cell (AC1) so the
# new header cells (AD1:AF1) get the same bold/centered/bordered style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# New header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (same for every player row) for rows 2-43.
$wins = 93
$losses = 69
$ties = 0

for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
